$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph "Play Bigger Bass Blizzard - Christmas Catch for
#    Free".
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Bigger Bass Blizzard – Christmas Catch and play for free. Discover its stunning graphics, high maximum win, and special Free Spins feature.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2. Remove the duplicated "Play Bigger Bass Blizzard - Christmas Catch for
#    Free" paragraph that used to sit right before the closing meta-text
#    paragraph near the end of the document.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$playAgainPara = $d.Paragraphs.Item($count - 1)
if ($playAgainPara.Range.Text -match "Play Bigger Bass Blizzard") {
    $playAgainPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3. Replace the final italic paragraph's text (previously the meta
#    description) with the new AI image-generation prompt, keeping the
#    italic formatting intact.
# ---------------------------------------------------------------------------
$oldText = "Read our review of Bigger Bass Blizzard – Christmas Catch and play for free. Discover its stunning graphics, high maximum win, and special Free Spins feature."
$newText = "Create an image featuring a happy Maya warrior with glasses in a cartoon style. The warrior should have a festive look, with a Santa hat, a scarf, and a fishing rod in hand, ready to catch some big bass in the frozen lake. In the background, there should be snow-covered trees and white flakes falling, creating a perfect Christmas atmosphere. The image should have bright and colorful tones to make it eye-catching and appealing to the players. The goal is to showcase the fun and thrilling experience of the game, while also highlighting the festive season and the unique character of the Maya warrior."

$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
